# Fruta / hortaliza, semanal
# Insert two new weekly records at rows 442-443 (pushing the existing
# rows 442-501 down to 444-503), then populate the two new rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before the current row 442; this shifts
# rows 442:501 down to 444:503 and extends the used range to T503.
$ws.Rows("442:443").Insert()

# --- New row 442 ---
$ws.Range("A442").Value = 9
$ws.Range("B442").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C442").Value = "Metropolitana"
$ws.Range("D442").Value = 44449
$ws.Range("E442").Value = 13
$ws.Range("F442").Value = "Fruta"
$ws.Range("G442").Value = 100104
$ws.Range("H442").Value = "Frutos de pepita"
$ws.Range("I442").Value = 100104005
$ws.Range("J442").Value = "Pera"
$ws.Range("K442").Value = "Packham's Triumph"
$ws.Range("L442").Value = "Especial"
$ws.Range("M442").Value = 580
$ws.Range("N442").Value = 14000
$ws.Range("O442").Value = 14000
$ws.Range("P442").Value = 14000
$ws.Range("Q442").Value = "$/caja 18 kilos granel"
$ws.Range("R442").Value = "Región de O'Higgins"
$ws.Range("S442").Value = 778
$ws.Range("T442").Value = 18

# --- New row 443 ---
$ws.Range("A443").Value = 9
$ws.Range("B443").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C443").Value = "Metropolitana"
$ws.Range("D443").Value = 44449
$ws.Range("E443").Value = 13
$ws.Range("F443").Value = "Fruta"
$ws.Range("G443").Value = 100104
$ws.Range("H443").Value = "Frutos de pepita"
$ws.Range("I443").Value = 100104005
$ws.Range("J443").Value = "Pera"
$ws.Range("K443").Value = "Packham's Triumph"
$ws.Range("L443").Value = "Primera"
$ws.Range("M443").Value = 450
$ws.Range("N443").Value = 12000
$ws.Range("O443").Value = 12000
$ws.Range("P443").Value = 12000
$ws.Range("Q443").Value = "$/caja 18 kilos granel"
$ws.Range("R443").Value = "Región de O'Higgins"
$ws.Range("S443").Value = 667
$ws.Range("T443").Value = 18
